$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sampling_Plan")
$ws.Range("C4").Value = "TestVal"
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats = -4122
Write-Host ("C4 value after format paste: " + $ws.Range("C4").Value2)
